$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level calculation settings: drop R1C1 refMode (calcPr refMode) ---
$excel.ReferenceStyle = 1  # xlA1

# --- Column widths for the 10 new vaccine columns (CY..DH / 103..112) ---
# Widths chosen so the engine-stored XML "width" lands as close as possible
# to the authored bestFit widths (20.332.., 31, 25.664.., 21.164.., 14.5).
$ws.Columns.Item(103).ColumnWidth = 19.571428571428573
$ws.Columns.Item(104).ColumnWidth = 30.285714285714285
$ws.Columns.Item(105).ColumnWidth = 25.0
$ws.Columns.Item(106).ColumnWidth = 20.428571428571427
$ws.Columns.Item(107).ColumnWidth = 13.714285714285714
$ws.Columns.Item(108).ColumnWidth = 19.571428571428573
$ws.Columns.Item(109).ColumnWidth = 30.285714285714285
$ws.Columns.Item(110).ColumnWidth = 25.0
$ws.Columns.Item(111).ColumnWidth = 20.428571428571427
$ws.Columns.Item(112).ColumnWidth = 13.714285714285714

# --- Vaccine Administration Date columns (DA, DF) are stored as Text ("@") ---
# Only format the rows that will actually receive data (1-7) so we don't
# materialize empty styled cells in rows 8-12, which the source file does not have.
$ws.Range("DA1:DA7").NumberFormat = "@"
$ws.Range("DF1:DF4").NumberFormat = "@"

# --- Header row + vaccine data for rows 1-7 ---
$ws.Range("CY1").Value = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value = "Vaccine 1 Product Name"
$ws.Range("DA1").Value = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value = "Vaccine 1 Dose Number"
$ws.Range("DC1").Value = "Vaccine 1 Notes"
$ws.Range("DD1").Value = "Vaccine 2 Group Name"
$ws.Range("DE1").Value = "Vaccine 2 Product Name"
$ws.Range("DF1").Value = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value = "Vaccine 2 Notes"

$ws.Range("CY2").Value = "COVID-19"
$ws.Range("CZ2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA2").Value = "2020-06-01"
$ws.Range("DB2").Value = 1
$ws.Range("DC2").Value = "notes 1"
$ws.Range("DD2").Value = "COVID-19"
$ws.Range("DE2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DF2").Value = "2020-06-20"
$ws.Range("DG2").Value = 2
$ws.Range("DH2").Value = "notes 2"

$ws.Range("CY3").Value = "COVID-19"
$ws.Range("CZ3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DA3").Value = "2020-06-02"
$ws.Range("DB3").Value = 1
$ws.Range("DD3").Value = "COVID-19"
$ws.Range("DE3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DF3").Value = "2020-06-21"
$ws.Range("DG3").Value = 2

$ws.Range("CY4").Value = "COVID-19"
$ws.Range("CZ4").Value = "Unknown"
$ws.Range("DA4").Value = "2020-06-04"
$ws.Range("DB4").Value = 1
$ws.Range("DD4").Value = "COVID-19"
$ws.Range("DE4").Value = "Unknown"
$ws.Range("DF4").Value = "2020-06-22"
$ws.Range("DG4").Value = 2

$ws.Range("CY5").Value = "COVID-19"
$ws.Range("CZ5").Value = "Moderna COVID-19 Vaccine"
$ws.Range("DA5").Value = "2020-06-01"
$ws.Range("DB5").Value = 1

$ws.Range("CY6").Value = "COVID-19"
$ws.Range("CZ6").Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("DA6").Value = "2020-06-03"
$ws.Range("DB6").Value = 1

$ws.Range("CY7").Value = "COVID-19"
$ws.Range("CZ7").Value = "Unknown"
$ws.Range("DA7").Value = "2020-06-02"
$ws.Range("DB7").Value = 1

# --- Reset the view: scroll back to top-left and select A1 (closest achievable
#     match for the saved file no longer being scrolled/selected deep in the sheet) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
